# Removing less than USD 5 price from extrapolation calibration because it is
# just a noise. This updates the recomputed extrapolation outputs (columns
# D..H) for the affected rows in the calibration results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = @(119494.7520288914, -0.03204552360320836, 0.2545398186810214, -3.019765692185008, 29.31631896758674)
    7  = @(121432.7456996691, -0.03302672754561237, 0.2347164006722956, -0.7972459370300036, 6.055105716062216)
    8  = @(121924.6873696327, -0.04917568539811389, 0.2080523357338072, -0.7777516457370579, 6.578587421454435)
    9  = @(124086.7039635979, -0.07663443972923126, 0.3374900485706538, -1.592942164420613, 10.32990551740569)
    10 = @(125316.744510438,  -0.113240112366871,   0.4305583691303654, -1.897574950301293, 9.820113034899308)
    11 = @(127039.6284207866, -0.1901139740748309,  0.752852153138254,  -2.605715163127442, 12.73584718709789)
    13 = @(117239.6927717878, -0.04358959310176971, 0.1312340009540396, -0.8517391197987699, 7.967351353277716)
    14 = @(117165.1276435176, -0.063930288213675,   0.1157528327797027, -1.350851472204498, 8.931316976582947)
    17 = @(117264.8404666179, -0.04708520277054174, 0.1380061068974084, -0.9571117429892404, 9.835338132462677)
    20 = @(118126.2829171024, -0.01615442316895325, 0.1467957770552073, -0.4985937497794729, 6.61528812132168)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
    $ws.Cells.Item($row, 6).Value = $vals[2]
    $ws.Cells.Item($row, 7).Value = $vals[3]
    $ws.Cells.Item($row, 8).Value = $vals[4]
}
